$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "user"
$ws.Range("B1").Value = "assistant"
$ws.Range("A2").Value = "I can write a for loop"
$ws.Range("B2").Value = "Great, Write a loop to display  a range of numbers from 1 to 10"
$ws.Range("B3").Value = "OK let’s start by learning the building blocks of the for loop. Lets start by looking at variables, define a variable that will store the number 5?"
$ws.Range("B4").Value = "You are right, thats exactly what loops are. To use loops we need to know the loop syntax, write for me the for loop syntax?"
$ws.Range("B6").Value = "You are right, write for me the syntax of a for loop`n"

$ws.Range("B9").Select()
$ws.Range("B9").Activate()
